# SMSTEST.xlsx - "We b issues fixed"
# Updates the Otcandnonrx sheet: A4/A5 get new text values (and B4/B5 follow),
# and the window selection moves to B11; also restores/updates the window size.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Otcandnonrx")
$ws.Activate()

$ws.Range("A4").Value = "2mg Tablet 10'S"
$ws.Range("B4").Value = "Gemer 2mg Tablet 10'S"
$ws.Range("A5").Value = "Dolo 650mg Tablet 15 Gemer "
$ws.Range("B5").Value = "Dolo 650mg Tablet 15 Gemer "

$ws.Range("B11").Select()

$excel.Width = 23256
$excel.Height = 12576
$excel.Left = -108
$excel.Top = -108
